# removed shopper role from test data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 holds the "AC_shopper " / "aa000g0l" / "Shopper" / 14 record.
# Deleting the entire row shifts the following row (AC_observer / ...) up
# so it becomes row 11, matching the target data.
$ws.Rows.Item(11).Delete()
